$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: target cluster changes from FAPs to ECs, plus several numeric updates ---
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.538334
$ws.Range("H2").Value = 7.615002
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.708201
$ws.Range("N2").Value = 2.124603
$ws.Range("O2").Value = 0.03793614316565257
$ws.Range("P2").Value = 0.03793614316565257
$ws.Range("Q2").Value = 1.797650677134
$ws.Range("R2").Value = 16.178856094206
$ws.Range("S2").Value = 0.03793614316565257
$ws.Range("T2").Value = 0.03793614316565257

# --- Row 3: now targets FAPs (previously targeted sCs) and gets new numbers ---
$ws.Range("A3").Value = "sCs"
$ws.Range("B3").Value = "Gdnf"
$ws.Range("C3").Value = "Gfra2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.538334
$ws.Range("H3").Value = 7.615002
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 13.05272433333334
$ws.Range("N3").Value = 39.15817300000001
$ws.Range("O3").Value = 0.699194182175866
$ws.Range("P3").Value = 0.6991941821758659
$ws.Range("Q3").Value = 33.13217396792734
$ws.Range("R3").Value = 298.189565711346
$ws.Range("S3").Value = 0.699194182175866
$ws.Range("T3").Value = 0.6991941821758659

# --- Row 4 (new): targets sCs ---
$ws.Range("A4").Value = "sCs"
$ws.Range("B4").Value = "Gdnf"
$ws.Range("C4").Value = "Gfra2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.538334
$ws.Range("H4").Value = 7.615002
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.907314
$ws.Range("N4").Value = 14.721942
$ws.Range("O4").Value = 0.2628696746584814
$ws.Range("P4").Value = 0.2628696746584814
$ws.Range("Q4").Value = 12.456401974876
$ws.Range("R4").Value = 112.107617773884
$ws.Range("S4").Value = 0.2628696746584814
$ws.Range("T4").Value = 0.2628696746584814
